$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness (column C) values for generations 0..250 (rows 2..252)
$newValues = @(12691,9919,9919,9919,9224,9209,9189,9189,8706,8706,8706,8706,8706,8706,8706,8706,8706,8706,8706,8706,8655,8655,8569,8569,8569,8569,8560,8560,8560,8233,8233,8233,8233,8233,7917,7892,7892,7892,7892,7892,7892,7892,7892,7892,7892,7892,7892,7892,7860,7860,7860,7860,7860,7750,7750,7750,7750,7750,7750,7750,7750,7750,7707,7707,7623,7623,7623,7623,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
